$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 97, pushing the existing weekly records
# (rows 97-144) down to rows 98-145, making room for the newest
# "Ajo" price-report entry for Terminal La Palmera de La Serena.
$ws.Rows.Item(97).EntireRow.Insert()

$ws.Range("A97").Value = 8
$ws.Range("B97").Value = "Terminal La Palmera de La Serena"
$ws.Range("C97").Value = "Coquimbo"
$ws.Range("D97").Value = 44460
$ws.Range("E97").Value = 4
$ws.Range("F97").Value = 100112003
$ws.Range("G97").Value = "Ajo"
$ws.Range("H97").Value = "Chino"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 540
$ws.Range("K97").Value = 15000
$ws.Range("L97").Value = 16000
$ws.Range("M97").Value = 15500
$ws.Range("N97").Value = "`$/caja 10 kilos"
$ws.Range("O97").Value = "China"
$ws.Range("P97").Value = 1550
$ws.Range("Q97").Value = 10
$ws.Range("R97").Value = "Hortaliza"
